$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "09.01.2026 12:45"
$ws.Range("D4").Value = "Termine"
